# X2A: introduce keypad — bump the char-count input (D5) from 11 to 12,
# clear the stray "cc" label in D6, and move the selection to the
# recalculated result cell D9.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 12
$ws.Range("D6").Value = $null

$ws.Range("D9").Select()
